$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 1186.8
$ws.Range("I19").Value = 547.25
$ws.Range("K19").Value = 547.25
$ws.Range("M19").Value = -372.25
# row 33
$ws.Range("H33").Value = 76924160
$ws.Range("I33").Value = 903.1
$ws.Range("K33").Value = 903.1
$ws.Range("M33").Value = -674.1
# row 38
$ws.Range("H38").Value = 6632.2
$ws.Range("I38").Value = 1264.7142
$ws.Range("K38").Value = 3794.1426
$ws.Range("M38").Value = -3422.1426
# row 39
$ws.Range("H39").Value = 314.5
$ws.Range("I39").Value = 277.6
$ws.Range("J39").Value = 499
$ws.Range("K39").Value = 832.8000000000001
$ws.Range("L39").Value = 1497
$ws.Range("M39").Value = -536.8000000000001
$ws.Range("N39").Value = -2089
# row 40
$ws.Range("H40").Value = 2486.375
$ws.Range("I40").Value = 2486.375
$ws.Range("K40").Value = 2486.375
$ws.Range("M40").Value = -2311.375
# row 42
$ws.Range("H42").Value = 2695.3333
$ws.Range("I42").Value = 27.666666
$ws.Range("K42").Value = 82.99999800000001
$ws.Range("M42").Value = 147.000002
# row 43
$ws.Range("H43").Value = 2795.7144
$ws.Range("I43").Value = 1892.75
$ws.Range("J43").Value = 3999.6667
$ws.Range("K43").Value = 1892.75
$ws.Range("L43").Value = 3999.6667
$ws.Range("M43").Value = -1823.75
$ws.Range("N43").Value = -4137.6667
# row 45
$ws.Range("H45").Value = 4800
$ws.Range("J45").Value = 4800
$ws.Range("L45").Value = 14400
$ws.Range("N45").Value = -14784
# row 49
$ws.Range("H49").Value = 1066.6666
$ws.Range("J49").Value = 500
$ws.Range("L49").Value = 1500
$ws.Range("N49").Value = -1772
# row 52
$ws.Range("H52").Value = 3909.2
$ws.Range("I52").Value = 4274
$ws.Range("J52").Value = 3666
$ws.Range("K52").Value = 12822
$ws.Range("L52").Value = 10998
$ws.Range("M52").Value = -12662
$ws.Range("N52").Value = -11318
# row 55
$ws.Range("H55").Value = 543.0833
$ws.Range("J55").Value = 909.8333
$ws.Range("L55").Value = 909.8333
$ws.Range("N55").Value = -1337.8333
# row 80
$ws.Range("H80").Value = 274.07693
$ws.Range("J80").Value = 252.4
$ws.Range("L80").Value = 757.2
$ws.Range("N80").Value = -2753.2
# row 83
$ws.Range("H83").Value = 274.07693
$ws.Range("J83").Value = 252.4
$ws.Range("L83").Value = 2271.6
$ws.Range("N83").Value = -12255.6
# row 132
$ws.Range("H132").Value = 2116
$ws.Range("I132").Value = 1976.0667
$ws.Range("K132").Value = 5928.2001
$ws.Range("M132").Value = -3398.2001
# row 141
$ws.Range("H141").Value = 3132.0527
$ws.Range("I141").Value = 2250.5
$ws.Range("K141").Value = 6751.5
$ws.Range("M141").Value = -1571.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 2554.9834
$ws.Range("I32").Value = 1011.86957
$ws.Range("K32").Value = 1011.86957
$ws.Range("M32").Value = -724.86957
# row 61
$ws.Range("H61").Value = 2341
$ws.Range("I61").Value = 1999
$ws.Range("J61").Value = 2512
$ws.Range("K61").Value = 1999
$ws.Range("L61").Value = 2512
$ws.Range("M61").Value = -1787
$ws.Range("N61").Value = -2936
# row 117
$ws.Range("H117").Value = 145000
$ws.Range("J117").Value = 145000
$ws.Range("L117").Value = 145000
$ws.Range("N117").Value = -154178
# row 136
$ws.Range("H136").Value = 2341
$ws.Range("I136").Value = 1999
$ws.Range("J136").Value = 2512
$ws.Range("K136").Value = 5997
$ws.Range("L136").Value = 7536
$ws.Range("M136").Value = -3447
$ws.Range("N136").Value = -12636

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 26
$ws.Range("H26").Value = 13956.667
$ws.Range("I26").Value = 13956.667
$ws.Range("K26").Value = 13956.667
$ws.Range("M26").Value = -13664.667
# row 86
$ws.Range("H86").Value = 2653.5
$ws.Range("I86").Value = 2479.5715
$ws.Range("K86").Value = 2479.5715
$ws.Range("M86").Value = -1356.5715
# row 89
$ws.Range("H89").Value = 2653.5
$ws.Range("I89").Value = 2479.5715
$ws.Range("K89").Value = 12397.8575
$ws.Range("M89").Value = -6781.8575

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
# row 7
$ws.Range("H7").Value = 75.28570999999999
$ws.Range("I7").Value = 53
$ws.Range("J7").Value = 92
$ws.Range("K7").Value = 53
$ws.Range("L7").Value = 92
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = -318
# row 134
$ws.Range("H134").Value = 1901.84
$ws.Range("I134").Value = 1708.7778
$ws.Range("K134").Value = 5126.3334
$ws.Range("M134").Value = -2591.3334

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 3
$ws.Range("H3").Value = 1437
$ws.Range("I3").Value = 1437
$ws.Range("K3").Value = 4311
$ws.Range("M3").Value = -4199
# row 33
$ws.Range("H33").Value = 390.53845
$ws.Range("I33").Value = 262.4
$ws.Range("J33").Value = 470.625
$ws.Range("K33").Value = 1574.4
$ws.Range("L33").Value = 2823.75
$ws.Range("M33").Value = -1291.4
$ws.Range("N33").Value = -3389.75
# row 36
$ws.Range("H36").Value = 1200
$ws.Range("I36").Value = 800
$ws.Range("K36").Value = 2400
$ws.Range("M36").Value = -2231
# row 68
$ws.Range("H68").Value = 998
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 998
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2994
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4616
# row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# row 71
$ws.Range("H71").Value = 998
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 998
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 8982
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -17094
# row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# row 131
$ws.Range("H131").Value = 10528055
$ws.Range("I131").Value = 7693234
$ws.Range("K131").Value = 23079702
$ws.Range("M131").Value = -23074662
# row 133
$ws.Range("H133").Value = 6600
$ws.Range("I133").Value = 8000
$ws.Range("K133").Value = 24000
$ws.Range("M133").Value = -18940
# row 140
$ws.Range("H140").Value = 858.7778
$ws.Range("I140").Value = 858.7778
$ws.Range("K140").Value = 2576.3334
$ws.Range("M140").Value = 2603.6666

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 15691.7
$ws.Range("I70").Value = 19737.742
$ws.Range("K70").Value = 19737.742
$ws.Range("M70").Value = -19467.742
# row 73
$ws.Range("H73").Value = 15691.7
$ws.Range("I73").Value = 19737.742
$ws.Range("K73").Value = 19737.742
$ws.Range("M73").Value = -18801.742
# row 107
$ws.Range("H107").Value = 963.2174
$ws.Range("I107").Value = 968.41174
$ws.Range("K107").Value = 968.41174
$ws.Range("M107").Value = 951.58826
# row 122
$ws.Range("H122").Value = 4332.8184
$ws.Range("I122").Value = 3844.25
$ws.Range("J122").Value = 5635.6665
$ws.Range("K122").Value = 11532.75
$ws.Range("L122").Value = 16906.9995
$ws.Range("M122").Value = -9082.75
$ws.Range("N122").Value = -21806.9995

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 1499.6666
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 1499
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 1499
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -1875
# row 68
$ws.Range("H68").Value = 2471.8572
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# row 71
$ws.Range("H71").Value = 2471.8572
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# row 122
$ws.Range("H122").Value = 6301.16
$ws.Range("I122").Value = 5797.3687
$ws.Range("K122").Value = 17392.1061
$ws.Range("M122").Value = -14942.1061

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 3017.1428
$ws.Range("I81").Value = 1182.7778
$ws.Range("J81").Value = 6319
$ws.Range("K81").Value = 2365.5556
$ws.Range("L81").Value = 12638
$ws.Range("M81").Value = -1304.5556
$ws.Range("N81").Value = -14760
# row 84
$ws.Range("H84").Value = 3017.1428
$ws.Range("I84").Value = 1182.7778
$ws.Range("J84").Value = 6319
$ws.Range("K84").Value = 11827.778
$ws.Range("L84").Value = 63190
$ws.Range("M84").Value = -6523.778
$ws.Range("N84").Value = -73798
# row 113
$ws.Range("H113").Value = 961.36365
$ws.Range("I113").Value = 1000.6875
$ws.Range("J113").Value = 856.5
$ws.Range("K113").Value = 3002.0625
$ws.Range("L113").Value = 2569.5
$ws.Range("M113").Value = -832.0625
$ws.Range("N113").Value = -6909.5
# row 132
$ws.Range("H132").Value = 1417.9474
$ws.Range("I132").Value = 1378.9412
$ws.Range("J132").Value = 1749.5
$ws.Range("K132").Value = 4136.8236
$ws.Range("L132").Value = 5248.5
$ws.Range("M132").Value = -1606.8236
$ws.Range("N132").Value = -10308.5
# row 136
$ws.Range("H136").Value = 2903.6843
$ws.Range("I136").Value = 2516.875
$ws.Range("K136").Value = 7550.625
$ws.Range("M136").Value = -5000.625
